$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Creating classes for Customers and its Manager
$ws.Range("B7").Value = "Thomas, Yamid"
$ws.Range("D7").Value = "In Progress"
$ws.Range("E7").Value = "Class Creation done (see src) and manager in progress"

# Row 8: Creating classes for Bills and its Manager(if required)
$ws.Range("B8").Value = "Thomas, Yamid"
$ws.Range("D8").Value = "In Progress"
$ws.Range("E8").Value = "Class Creation done (see src) and manager in progress (group with Customer)"

# Move the active selection to D8, matching the saved cursor position
$ws.Range("D8").Select()
